$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 74, shifting existing rows 74..146 down to 75..147
$ws.Rows.Item(74).Insert()

# Populate the newly inserted row 74 with its data
$ws.Range("A74").Value = 6
$ws.Range("B74").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C74").Value = "Metropolitana"
$ws.Range("D74").Value = 44484
$ws.Range("E74").Value = 13
$ws.Range("F74").Value = 100112022
$ws.Range("G74").Value = "Arveja Verde"
$ws.Range("H74").Value = "Perfection"
$ws.Range("I74").Value = "Primera"
$ws.Range("J74").Value = 280
$ws.Range("K74").Value = 18000
$ws.Range("L74").Value = 20000
$ws.Range("M74").Value = 18929
$ws.Range("N74").Value = "$/malla 25 kilos"
$ws.Range("O74").Value = "Provincia de Limarí"
$ws.Range("P74").Value = 757
$ws.Range("Q74").Value = 25
$ws.Range("R74").Value = "Hortaliza"
